$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values would otherwise be
# auto-detected as numbers by Excel (so formatting like trailing
# zeros / thousand-dot separators in the source data survives).
$textCells = @("D5", "D11", "D16", "D19", "D20", "D21", "D25", "D27", "D29", "D39", "D40", "D41", "D42", "D43", "D44", "D47", "D49", "D51")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.225.46"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "1.589.01"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "212.02"
$ws.Range("E5").Value = "  +0.85%  "
$ws.Range("E6").Value = "  -0.47%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -1.11%  "
$ws.Range("E10").Value = "  -1.82%  "
$ws.Range("D11").Value = "0.0846"
$ws.Range("E11").Value = "  +0.39%  "
$ws.Range("D12").Value = "1.813.16"
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("D13").Value = "1.617.92"
$ws.Range("E13").Value = "  +1.70%  "
$ws.Range("E14").Value = "  -1.41%  "
$ws.Range("E15").Value = "  -0.19%  "
$ws.Range("D16").Value = "64.02"
$ws.Range("E16").Value = "  -0.92%  "
$ws.Range("D17").Value = "26.235.04"
$ws.Range("E17").Value = "  -0.36%  "
$ws.Range("D18").Value = "0.0₃0724"
$ws.Range("E18").Value = "  -0.62%  "
$ws.Range("D19").Value = "214.63"
$ws.Range("E19").Value = "  +1.20%  "
$ws.Range("D20").Value = "7.28"
$ws.Range("E20").Value = "  -2.80%  "
$ws.Range("D21").Value = "1.00"
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("E22").Value = "  -0.84%  "
$ws.Range("E23").Value = "  -0.65%  "
$ws.Range("E24").Value = "  +0.48%  "
$ws.Range("D25").Value = "144.00"
$ws.Range("E25").Value = "  -0.86%  "
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("D27").Value = "7.00"
$ws.Range("E27").Value = "  -0.77%  "
$ws.Range("E28").Value = "  -0.81%  "
$ws.Range("D29").Value = "15.15"
$ws.Range("E29").Value = "  -0.99%  "
$ws.Range("E30").Value = "  -1.80%  "
$ws.Range("E31").Value = "  +0.72%  "
$ws.Range("E32").Value = "  -1.22%  "
$ws.Range("D33").Value = "1.386.15"
$ws.Range("E33").Value = "  +6.59%  "
$ws.Range("E34").Value = "  -1.67%  "
$ws.Range("E35").Value = "  -0.49%  "
$ws.Range("E36").Value = "  -1.23%  "
$ws.Range("E37").Value = "  -5.64%  "
$ws.Range("E38").Value = "  -0.91%  "
$ws.Range("D39").Value = "0.819"
$ws.Range("E39").Value = "  +0.65%  "
$ws.Range("D40").Value = "5.84"
$ws.Range("E40").Value = "  +4.01%  "
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "0.931"
$ws.Range("E42").Value = "  -15.71%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "0.767"
$ws.Range("E43").Value = "  +0.73%  "
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").Value = "2.14"
$ws.Range("E44").Value = "  -0.17%  "
$ws.Range("D45").Value = "1.725.25"
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("E46").Value = "  -2.58%  "
$ws.Range("D47").Value = "86.11"
$ws.Range("E47").Value = "  -2.68%  "
$ws.Range("E48").Value = "  -1.28%  "
$ws.Range("D49").Value = "0.0977"
$ws.Range("E49").Value = "  -1.03%  "
$ws.Range("E50").Value = "  -1.04%  "
$ws.Range("D51").Value = "0.999"
$ws.Range("E51").Value = "  -0.21%  "

Write-Output "Applied cryptos list update"
